# Update the power-loss results table (rows 2-25, columns B:N excluding
# F/H/K/O which remain 0) with the recomputed values for the 380 kV case.
$data = @{
    2 = @{ 'B'=5.261984186751761; 'C'=0.9858216608980683; 'D'=0.2914293965006038; 'E'=0.04087277362430974; 'G'=0.002743893940322697; 'I'=3.91773763586292; 'J'=0.02480685757290502; 'L'=0.767596745349735; 'M'=1.144258733681184; 'N'=4.082271666634526 }
    3 = @{ 'B'=5.171146224688584; 'C'=0.9506475095675455; 'D'=0.291997607752684; 'E'=0.04051295017199052; 'G'=0.002751313913116166; 'I'=3.901324815216753; 'J'=0.02163968271359806; 'L'=0.7649066336240793; 'M'=1.13017328850885; 'N'=4.08578474987327 }
    4 = @{ 'B'=5.118552558480701; 'C'=0.9297073661742843; 'D'=0.2924477854757797; 'E'=0.04028737989208775; 'G'=0.002756106599938953; 'I'=3.892793537868542; 'J'=0.01970034907523655; 'L'=0.7636045592358158; 'M'=1.122168892013157; 'N'=4.088780339964885 }
    5 = @{ 'B'=5.097917246920588; 'C'=0.9213378683924702; 'D'=0.2926567203510047; 'E'=0.04019427670967435; 'G'=0.002758119426760618; 'I'=3.889704448001808; 'J'=0.01891127154877381; 'L'=0.7631617448616907; 'M'=1.119068526344613; 'N'=4.090211076252714 }
    6 = @{ 'B'=5.094538824753556; 'C'=0.9199579709083423; 'D'=0.2926929533847868; 'E'=0.0401787451922706; 'G'=0.002758457271836725; 'I'=3.889214870722029; 'J'=0.01878031576568162; 'L'=0.7630935144580917; 'M'=1.118563452931824; 'N'=4.090461308285413 }
    7 = @{ 'B'=5.118271040048285; 'C'=0.9295938308993072; 'D'=0.2924505000413333; 'E'=0.0402861290739569; 'G'=0.002756133503444853; 'I'=3.892750310195481; 'J'=0.01968970253474822; 'L'=0.7635982319917218; 'M'=1.122126426144419; 'N'=4.088798786163593 }
    8 = @{ 'B'=5.230001182144406; 'C'=0.9735564002015167; 'D'=0.2916042990819321; 'E'=0.04074965958633392; 'G'=0.002746403335761199; 'I'=3.911756739131206; 'J'=0.02371363574963681; 'L'=0.7665965363945304; 'M'=1.139268084174795; 'N'=4.083308421011964 }
    9 = @{ 'B'=5.474498011439721; 'C'=1.065049268496352; 'D'=0.2907483097286203; 'E'=0.04162257137358161; 'G'=0.002729191106925277; 'I'=3.961364391245738; 'J'=0.03165272084328308; 'L'=0.7752581678543322; 'M'=1.178019131995399; 'N'=4.079235320739798 }
    10 = @{ 'B'=5.669854275597515; 'C'=1.135595586778834; 'D'=0.2906091413718315; 'E'=0.04224295621594454; 'G'=0.002717670050929315; 'I'=4.005434168442591; 'J'=0.03752370626780532; 'L'=0.7833303815525028; 'M'=1.209661643511424; 'N'=4.080381284670011 }
    11 = @{ 'B'=5.762197328310492; 'C'=1.168435194375036; 'D'=0.290652209566673; 'E'=0.04252086790011145; 'G'=0.002712669975162434; 'I'=4.027161220596554; 'J'=0.04020480505469948; 'L'=0.787376593807096; 'M'=1.224754892384709; 'N'=4.081814148644327 }
    12 = @{ 'B'=5.797669038642368; 'C'=1.180980050266612; 'D'=0.29068381625261; 'E'=0.04262550488304218; 'G'=0.002710810981408782; 'I'=4.035632001454474; 'J'=0.0412217070514842; 'L'=0.78896280829197; 'M'=1.230571519511201; 'N'=4.082488865908942 }
    13 = @{ 'B'=5.790007120982182; 'C'=1.178273406741482; 'D'=0.290676328787896; 'E'=0.0426029959464671; 'G'=0.002711209821292655; 'I'=4.033796822368529; 'J'=0.04100262437127355; 'L'=0.7886187835486282; 'M'=1.22931429609892; 'N'=4.082337660467545 }
    14 = @{ 'B'=5.765105496364754; 'C'=1.169465068711986; 'D'=0.2906545032678736; 'E'=0.04252948841804738; 'G'=0.002712516345702704; 'I'=4.027853232450653; 'J'=0.04028843276270067; 'L'=0.787506009156445; 'M'=1.225231399042372; 'N'=4.081867004071995 }
    15 = @{ 'B'=5.749918203743391; 'C'=1.164083980279315; 'D'=0.2906431267721103; 'E'=0.04248438502176644; 'G'=0.002713321106797977; 'I'=4.024244337587476; 'J'=0.03985118595731763; 'L'=0.7868314415146784; 'M'=1.222743697609161; 'N'=4.081595950956199 }
    16 = @{ 'B'=5.663889675562359; 'C'=1.133464630779599; 'D'=0.290608467371456; 'E'=0.04222470885077012; 'G'=0.00271800164691582; 'I'=4.004048195241225; 'J'=0.03734871038265908; 'L'=0.7830734975480311; 'M'=1.208689382207353; 'N'=4.080306075148741 }
    17 = @{ 'B'=5.612006510805998; 'C'=1.114873352034692; 'D'=0.2906144506611525; 'E'=0.04206431535881183; 'G'=0.002720934555788163; 'I'=3.992089876250859; 'J'=0.03581627350767036; 'L'=0.7808640842466872; 'M'=1.200246970756268; 'N'=4.07974898542966 }
    18 = @{ 'B'=5.582491491389192; 'C'=1.104250417373066; 'D'=0.2906279046098774; 'E'=0.04197165508457879; 'G'=0.002722644176095718; 'I'=3.985369744207489; 'J'=0.0349358228709491; 'L'=0.7796284895764103; 'M'=1.195456864394046; 'N'=4.079514299028205 }
    19 = @{ 'B'=5.572554228944341; 'C'=1.100665693628514; 'D'=0.2906341795447531; 'E'=0.04194021163683637; 'G'=0.002723226927392347; 'I'=3.983121504204348; 'J'=0.03463787923044492; 'L'=0.7792161782683422; 'M'=1.193846289397428; 'N'=4.079449532172589 }
    20 = @{ 'B'=5.617495712289156; 'C'=1.116845139883139; 'D'=0.2906127774712246; 'E'=0.04208143143668419; 'G'=0.002720619995938164; 'I'=3.993346497488346; 'J'=0.03597930300104935; 'L'=0.7810956354055492; 'M'=1.201138870909588; 'N'=4.079799407660374 }
    21 = @{ 'B'=5.772406014364947; 'C'=1.172049313579691; 'D'=0.2906604987623922; 'E'=0.04255109556757386; 'G'=0.002712131655268118; 'I'=4.029592394186409; 'J'=0.04049816284921803; 'L'=0.7878313907137766; 'M'=1.226427894735863; 'N'=4.082001653089378 }
    22 = @{ 'B'=5.87658488812292; 'C'=1.208765796600858; 'D'=0.2907808556189195; 'E'=0.04285454647309361; 'G'=0.002706784598590607; 'I'=4.05469966607069; 'J'=0.04346105465083383; 'L'=0.792548460172938; 'M'=1.243545574490369; 'N'=4.084211557517818 }
    23 = @{ 'B'=5.820712738645113; 'C'=1.189110632838322; 'D'=0.2907084592918352; 'E'=0.04269290380108171; 'G'=0.002709620143060388; 'I'=4.041169075259987; 'J'=0.04187878278575852; 'L'=0.7900019918547088; 'M'=1.234355359595085; 'N'=4.082961229876389 }
    24 = @{ 'B'=5.615013068712813; 'C'=1.115953490895436; 'D'=0.2906135027287391; 'E'=0.04207369465429966; 'G'=0.002720762135549433; 'I'=3.992777896786265; 'J'=0.0359055956092007; 'L'=0.7809908432877535; 'M'=1.200735444959065; 'N'=4.079776345248973 }
    25 = @{ 'B'=5.405611735589105; 'C'=1.039721650326157; 'D'=0.2908938780362078; 'E'=0.04139017977577542; 'G'=0.002733648912774312; 'I'=3.946613077186583; 'J'=0.02949893194779207; 'L'=0.7726158701137251; 'M'=1.166981914426728; 'N'=4.079614691953026 }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
